# The diary/log entry document ends with a paragraph about the LAN
# simulation ("Queria intentar simular..."). That paragraph currently
# carries the document's "_GoBack" bookmark at its very end (right after
# the paragraph's text run).
#
# The edit appends a new log entry to the document:
#   - a blank paragraph
#   - a paragraph containing "10/9 hora entrada: 13.00"
# and the "_GoBack" bookmark moves along to sit at the new end of the
# document (inside the newly added last paragraph).

$d = $word.ActiveDocument

# Remove the existing "_GoBack" bookmark - it will be recreated at the
# new end of the document as part of the inserted content below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Collapsed range positioned right before the final paragraph mark of the
# document (i.e. at the very end of the last paragraph's text).
$endPos = $d.Content.End
$tail = $d.Range($endPos - 1, $endPos - 1)

# Build the two new paragraphs as raw OOXML (matching the formatting used
# throughout the rest of the document) and insert them in one shot. Using
# InsertXML lets us produce a genuinely empty paragraph (no run element)
# for the blank line, and keeps the new run/bookmark ordering exact.
$newParagraphs = '<w:p><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="es-ES"/></w:rPr></w:pPr></w:p>' +
                  '<w:p><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="es-ES"/></w:rPr></w:pPr>' +
                  '<w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="es-ES"/></w:rPr>' +
                  '<w:t>10/9 hora entrada: 13.00</w:t></w:r>' +
                  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
              '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
              '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
              '<w:body>' + $newParagraphs + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$null = $tail.InsertXML($packageXml)
